$wb = $excel.ActiveWorkbook

# --- Locate the existing "0.10.0" sheet (the newest protocol sheet so far) ---
$src = $wb.Worksheets.Item("0.10.0")

# --- Duplicate it, placing the copy right after it, then rename to "0.10.5" ---
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "0.10.5"

# --- Reset the frozen-pane scroll position back to the top of the data (B3) ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 2

# --- Row 9 (test #7 "Selecting an Artifact"): now fails ---
$ws.Range("C9").Value = "fail"
$ws.Range("C9").Interior.Color = 255
$ws.Range("D9").Value = "Artifact selected, but editor is not set dirty"

# --- Row 11 (test #9 "Adding a new Type..."): comment added, highlight changed to yellow ---
$ws.Range("C11").Interior.Color = 65535
$ws.Range("D11").Value = "HardwareAlgorithm and SubPipelineAlgorithm not tested"

# --- Row 12 (test #10 "Configuring a new Type..."): clear the old highlight/comment ---
$ws.Range("D21").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D12").ClearContents()

# --- Row 13 (test #11 "Deleting a Type..."): now fails, new highlight + comment ---
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "fail"
$ws.Range("C13").Interior.Color = 255
$ws.Range("D13").Value = "For Algorithms error- eu/qualimaster/observables/Imeasurable after selecting class; same for Sink and Source"

# --- Row 15 (test #13 "Adding a pipeline"): now fails ---
$ws.Range("C15").Value = "fail"
$ws.Range("C15").Interior.Color = 255
$ws.Range("D15").Value = "was not able to delete FamilyElement and algorithm, after restart app crashed with unknown elements"

# --- Row 16 (test #14 "Deleting a pipeline"): highlighted yellow, comment added ---
$ws.Range("C16").Interior.Color = 65535
$ws.Range("D16").Value = "It is possible to clone FamilyElement in Algoprithms, and it appears only in Algorithms, not under Families"

# --- Row 21 (test #19 "Validating the model"): now fails, comment replaced, row height reset ---
$ws.Range("C21").Value = "fail"
$ws.Range("C21").Interior.Color = 255
$ws.Range("D21").Value = 'adding name == "roman" eror name is unknow, same for latency == 1'
$ws.Rows.Item(21).AutoFit()

# --- Update the banner/title last (matches the order strings were appended upstream) ---
$ws.Range("C1").Value = "Roman: 0.10.5 win 64-bit (2016-07-01), Windows 7 64-bit, Java jdk1.8.0_74 64-bit with ConfModel"

# --- Make the new sheet the active / visible tab ---
$ws.Activate()
